# Roles workbook update:
#  - BranchUser role ("BranchUser" / "Standard Branch User", row 5) loses the
#    assetCreation and assetModification permissions (changed from TRUE to
#    FALSE) as part of tightening user-specific access for the new
#    notification system.
#  - Leave the cursor/selection on G11, matching the state the sheet was
#    saved in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Roles")

# Revoke assetCreation (D) and assetModification (E) permissions for the
# BranchUser role (row 5). Copy from a cell that already holds the text
# "false" so the written cells stay text (shared-string) values instead of
# being auto-converted to native Excel booleans.
$ws.Range("C5").Copy($ws.Range("D5"))
$ws.Range("C5").Copy($ws.Range("E5"))

# Restore the active selection/cursor position.
$ws.Range("G11").Select()
